$wb = $excel.ActiveWorkbook

# Update "展览" sheet: F2 902 -> 906, F3 68 -> 69
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 906
$ws1.Range("F3").Value = 69

# Update "全部类型" sheet: F2 902 -> 906, F3 68 -> 69
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 906
$ws2.Range("F3").Value = 69
